$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22: Topic 9 checklist (Feb 12) - add date (same date-number-format as D5:D21)
$ws.Range("D21").Copy()
$ws.Range("D22").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D22").Value = 45333

# Row 23: Topic 10 checklist (Feb 26) - add date, weight, and formula result changes
$ws.Range("D21").Copy()
$ws.Range("D23").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D23").Value = 45333
$ws.Range("F23").Value = 100

$excel.CutCopyMode = $false

# update selection to M17
$ws.Range("M17").Select()
